$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value2 = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "261.34"
Set-TextValue "D4" "6.207"
Set-TextValue "D5" "0.06163"
Set-TextValue "D6" "6.726"
Set-TextValue "D7" "3.464"
Set-TextValue "D8" "1.347"
Set-TextValue "D9" "0.7986"
Set-TextValue "B10" "One"
Set-TextValue "C10" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D10" "0.0006140"
Set-TextValue "E10" "9OneONEWorstin24h"
Set-TextValue "B11" "WazirX"
Set-TextValue "C11" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1586"
Set-TextValue "E11" "10WazirXWRX"
Set-TextValue "B12" "MandalaExchangeToken"
Set-TextValue "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.08113"
Set-TextValue "E12" "11MandalaExchangeTokenMDX"
Set-TextValue "B13" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C13" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D13" "0.03484"
Set-TextValue "E13" "12LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue "B14" "BitrueCoin"
Set-TextValue "C14" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D14" "0.03085"
Set-TextValue "E14" "13BitrueCoinBTR"
Set-TextValue "B15" "BitMartToken"
Set-TextValue "C15" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D15" "0.09308"
Set-TextValue "E15" "14BitMartTokenBMX"
Set-TextValue "B16" "MCDex"
Set-TextValue "C16" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D16" "3.848"
Set-TextValue "E16" "15MCDexMCB"
Set-TextValue "B17" "BitForexToken"
Set-TextValue "C17" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D17" "0.001688"
Set-TextValue "E17" "16BitForexTokenBF"
Set-TextValue "B18" "CoinExToken"
Set-TextValue "C18" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D18" "0.04793"
Set-TextValue "E18" "17CoinExTokenCET"
Set-TextValue "D19" "0.006192"
Set-TextValue "D20" "0.001092"
Set-TextValue "D21" "0.004071"
Set-TextValue "D22" "0.0001500"
Set-TextValue "D24" "2.207"
Set-TextValue "D27" "0.0003202"
Set-TextValue "D40" "0.04614"
Set-TextValue "D41" "0.007090"
Set-TextValue "D43" "0.003600"
Set-TextValue "D45" "0.002970"
Set-TextValue "D46" "0.00005937"
Set-TextValue "D48" "0.7000"
Set-TextValue "D49" "0.08931"
Set-TextValue "D50" "0.00002100"
Set-TextValue "D51" "0.01010"

Write-Output "Applied 58 cell updates"
